# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the d7a0cf63-... file (row 3) on the zh-cn and de-de
# language sheets, plus the corresponding "Latest HO Xliff Generate Date"
# roll-up on the Overview sheet, to reflect a freshly generated handback
# report.

$wb = $excel.ActiveWorkbook
$dateFormat = "yyyy-mm-dd HH:mm:ss"

# Overview sheet: Latest HO Xliff Generate Date for the d7a0cf63 row.
$wsOverview = $wb.Worksheets.Item("Overview")
$cellOverviewG3 = $wsOverview.Range("G3")
$cellOverviewG3.Value = "2016-08-16 08:48:18"
$cellOverviewG3.NumberFormat = $dateFormat

# zh-cn sheet: Correspond Handoff / Handback datetimes for the d7a0cf63 row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$cellZhCnH3 = $wsZhCn.Range("H3")
$cellZhCnH3.Value = "2016-08-16 08:48:13"
$cellZhCnH3.NumberFormat = $dateFormat
$cellZhCnK3 = $wsZhCn.Range("K3")
$cellZhCnK3.Value = "2016-08-16 08:48:30"
$cellZhCnK3.NumberFormat = $dateFormat

# de-de sheet: Correspond Handoff / Handback datetimes for the d7a0cf63 row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$cellDeDeH3 = $wsDeDe.Range("H3")
$cellDeDeH3.Value = "2016-08-16 08:48:18"
$cellDeDeH3.NumberFormat = $dateFormat
$cellDeDeK3 = $wsDeDe.Range("K3")
$cellDeDeK3.Value = "2016-08-16 08:48:37"
$cellDeDeK3.NumberFormat = $dateFormat
